# Apply the "456a3b4" gh-pages data refresh to 广州-漫展信息.xlsx
# Updates "想去人数" (F column) counters across sheets, and inserts a
# duplicated "NIJISANJI EN" row into the "全部类型" sheet (row 25),
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# Sheet "展览" (Exhibitions) - F column counter updates
# ---------------------------------------------------------------
$wsExpo = $wb.Worksheets.Item("展览")
$wsExpo.Range("F2").Value  = 829
$wsExpo.Range("F3").Value  = 13645
$wsExpo.Range("F4").Value  = 13443
$wsExpo.Range("F7").Value  = 35
$wsExpo.Range("F8").Value  = 584
$wsExpo.Range("F11").Value = 40
$wsExpo.Range("F12").Value = 736
$wsExpo.Range("F13").Value = 2124
$wsExpo.Range("F14").Value = 64
$wsExpo.Range("F16").Value = 65
$wsExpo.Range("F20").Value = 362
$wsExpo.Range("F22").Value = 494
$wsExpo.Range("F24").Value = 63

# ---------------------------------------------------------------
# Sheet "演出" (Performances) - F column counter update
# ---------------------------------------------------------------
$wsShow = $wb.Worksheets.Item("演出")
$wsShow.Range("F7").Value = 1349

# ---------------------------------------------------------------
# Sheet "本地生活" (Local Life) - F column counter updates
# ---------------------------------------------------------------
$wsLocal = $wb.Worksheets.Item("本地生活")
$wsLocal.Range("F2").Value = 213
$wsLocal.Range("F3").Value = 91

# ---------------------------------------------------------------
# Sheet "全部类型" (All Types) - F column counter updates (rows 2-24)
# ---------------------------------------------------------------
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F2").Value  = 213
$wsAll.Range("F3").Value  = 829
$wsAll.Range("F4").Value  = 13645
$wsAll.Range("F5").Value  = 13443
$wsAll.Range("F8").Value  = 35
$wsAll.Range("F9").Value  = 584
$wsAll.Range("F12").Value = 40
$wsAll.Range("F13").Value = 736
$wsAll.Range("F16").Value = 2124
$wsAll.Range("F17").Value = 64
$wsAll.Range("F19").Value = 65

# Row 24 ("NIJISANJI EN") counter bumps 90 -> 91 in place ...
$wsAll.Range("F24").Value = 91

# ... and then a duplicate of that same row is inserted right after it
# (new row 25), shifting the rest of the table (old rows 25-37) down to
# become rows 26-38. Copying row 24 preserves all of its formatting.
$wsAll.Rows.Item(24).Copy()
$wsAll.Rows.Item(25).Insert()
$excel.CutCopyMode = $false

# Fix up the new row's running index (A25) and counter (F25); the rest
# of the copied cells (B25:E25, G25:I25) already hold the correct data.
$wsAll.Range("A25").Value = 24
$wsAll.Range("F25").Value = 91
